$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "BTC"
$ws.Cells.Item(2, 3).Value = "Bitcoin"
$ws.Cells.Item(2, 4).Value = 51872
$ws.Cells.Item(2, 5).Value = 1017998543299
$ws.Cells.Item(2, 6).Value = 21205363662
$ws.Cells.Item(2, 7).Value = 0.35419
$ws.Cells.Item(3, 2).Value = "ETH"
$ws.Cells.Item(3, 3).Value = "Ethereum"
$ws.Cells.Item(3, 4).Value = 2809.62
$ws.Cells.Item(3, 5).Value = 337413977458
$ws.Cells.Item(3, 6).Value = 19945797889
$ws.Cells.Item(3, 7).Value = 0.9633699999999999
$ws.Cells.Item(4, 2).Value = "USDT"
$ws.Cells.Item(4, 3).Value = "Tether"
$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(4, 5).Value = 97705040238
$ws.Cells.Item(4, 6).Value = 43428728969
$ws.Cells.Item(4, 7).Value = 0.03638
$ws.Cells.Item(5, 2).Value = "BNB"
$ws.Cells.Item(5, 3).Value = "BNB"
$ws.Cells.Item(5, 4).Value = 354.9
$ws.Cells.Item(5, 5).Value = 54576196035
$ws.Cells.Item(5, 6).Value = 1037820454
$ws.Cells.Item(5, 7).Value = -1.14481
$ws.Cells.Item(6, 2).Value = "SOL"
$ws.Cells.Item(6, 3).Value = "Solana"
$ws.Cells.Item(6, 4).Value = 111.98
$ws.Cells.Item(6, 5).Value = 49348754862
$ws.Cells.Item(6, 6).Value = 1608657980
$ws.Cells.Item(6, 7).Value = 2.39227
$ws.Cells.Item(7, 2).Value = "XRP"
$ws.Cells.Item(7, 3).Value = "XRP"
$ws.Cells.Item(7, 4).Value = 0.558534
$ws.Cells.Item(7, 5).Value = 30475429761
$ws.Cells.Item(7, 6).Value = 1138165010
$ws.Cells.Item(7, 7).Value = 0.91177
$ws.Cells.Item(8, 2).Value = "USDC"
$ws.Cells.Item(8, 3).Value = "USDC"
$ws.Cells.Item(8, 4).Value = 1.001
$ws.Cells.Item(8, 5).Value = 28109207081
$ws.Cells.Item(8, 6).Value = 4437214552
$ws.Cells.Item(8, 7).Value = 0.13931
$ws.Cells.Item(9, 2).Value = "STETH"
$ws.Cells.Item(9, 3).Value = "Lido Staked Ether"
$ws.Cells.Item(9, 4).Value = 2806.63
$ws.Cells.Item(9, 5).Value = 27450979155
$ws.Cells.Item(9, 6).Value = 7579947
$ws.Cells.Item(9, 7).Value = 0.97972
$ws.Cells.Item(10, 2).Value = "ADA"
$ws.Cells.Item(10, 3).Value = "Cardano"
$ws.Cells.Item(10, 4).Value = 0.6345190000000001
$ws.Cells.Item(10, 5).Value = 22318054875
$ws.Cells.Item(10, 6).Value = 686428302
$ws.Cells.Item(10, 7).Value = 8.289239999999999
$ws.Cells.Item(11, 2).Value = "AVAX"
$ws.Cells.Item(11, 3).Value = "Avalanche"
$ws.Cells.Item(11, 4).Value = 40.42
$ws.Cells.Item(11, 5).Value = 14832671767
$ws.Cells.Item(11, 6).Value = 497055411
$ws.Cells.Item(11, 7).Value = 1.75858
$ws.Cells.Item(12, 2).Value = "DOGE"
$ws.Cells.Item(12, 3).Value = "Dogecoin"
$ws.Cells.Item(12, 4).Value = 0.08404300000000001
$ws.Cells.Item(12, 5).Value = 12037681282
$ws.Cells.Item(12, 6).Value = 331455934
$ws.Cells.Item(12, 7).Value = -0.56474
$ws.Cells.Item(13, 2).Value = "TRX"
$ws.Cells.Item(13, 3).Value = "TRON"
$ws.Cells.Item(13, 4).Value = 0.13545
$ws.Cells.Item(13, 5).Value = 11925909458
$ws.Cells.Item(13, 6).Value = 436704782
$ws.Cells.Item(13, 7).Value = -2.11869
$ws.Cells.Item(14, 2).Value = "LINK"
$ws.Cells.Item(14, 3).Value = "Chainlink"
$ws.Cells.Item(14, 4).Value = 20.08
$ws.Cells.Item(14, 5).Value = 11778615412
$ws.Cells.Item(14, 6).Value = 480891566
$ws.Cells.Item(14, 7).Value = 2.10767
$ws.Cells.Item(15, 2).Value = "DOT"
$ws.Cells.Item(15, 3).Value = "Polkadot"
$ws.Cells.Item(15, 4).Value = 7.78
$ws.Cells.Item(15, 5).Value = 10364588487
$ws.Cells.Item(15, 6).Value = 225924736
$ws.Cells.Item(15, 7).Value = 2.33205
$ws.Cells.Item(16, 2).Value = "MATIC"
$ws.Cells.Item(16, 3).Value = "Polygon"
$ws.Cells.Item(16, 4).Value = 0.947923
$ws.Cells.Item(16, 5).Value = 8797705421
$ws.Cells.Item(16, 6).Value = 408236220
$ws.Cells.Item(16, 7).Value = 1.87194
$ws.Cells.Item(17, 2).Value = "WBTC"
$ws.Cells.Item(17, 3).Value = "Wrapped Bitcoin"
$ws.Cells.Item(17, 4).Value = 51813
$ws.Cells.Item(17, 5).Value = 8126466520
$ws.Cells.Item(17, 6).Value = 114251529
$ws.Cells.Item(17, 7).Value = 0.3653
$ws.Cells.Item(18, 2).Value = "TON"
$ws.Cells.Item(18, 3).Value = "Toncoin"
$ws.Cells.Item(18, 4).Value = 2.25
$ws.Cells.Item(18, 5).Value = 7806563017
$ws.Cells.Item(18, 6).Value = 19313971
$ws.Cells.Item(18, 7).Value = 2.98088
$ws.Cells.Item(19, 2).Value = "ICP"
$ws.Cells.Item(19, 3).Value = "Internet Computer"
$ws.Cells.Item(19, 4).Value = 13.72
$ws.Cells.Item(19, 5).Value = 6296033599
$ws.Cells.Item(19, 6).Value = 107697648
$ws.Cells.Item(19, 7).Value = 3.74993
$ws.Cells.Item(20, 2).Value = "UNI"
$ws.Cells.Item(20, 3).Value = "Uniswap"
$ws.Cells.Item(20, 4).Value = 7.67
$ws.Cells.Item(20, 5).Value = 5789890433
$ws.Cells.Item(20, 6).Value = 176176979
$ws.Cells.Item(20, 7).Value = 0.93396
$ws.Cells.Item(21, 2).Value = "SHIB"
$ws.Cells.Item(21, 3).Value = "Shiba Inu"
$ws.Cells.Item(21, 4).Value = 0.00000977
$ws.Cells.Item(21, 5).Value = 5759077663
$ws.Cells.Item(21, 6).Value = 155364832
$ws.Cells.Item(21, 7).Value = 0.9297299999999999
$ws.Cells.Item(22, 2).Value = "BCH"
$ws.Cells.Item(22, 3).Value = "Bitcoin Cash"
$ws.Cells.Item(22, 4).Value = 267.59
$ws.Cells.Item(22, 5).Value = 5260435818
$ws.Cells.Item(22, 6).Value = 172266661
$ws.Cells.Item(22, 7).Value = -0.11901
$ws.Cells.Item(23, 2).Value = "LTC"
$ws.Cells.Item(23, 3).Value = "Litecoin"
$ws.Cells.Item(23, 4).Value = 70.39
$ws.Cells.Item(23, 5).Value = 5225289510
$ws.Cells.Item(23, 6).Value = 261464128
$ws.Cells.Item(23, 7).Value = 0.33215
$ws.Cells.Item(24, 2).Value = "DAI"
$ws.Cells.Item(24, 3).Value = "Dai"
$ws.Cells.Item(24, 4).Value = 0.999741
$ws.Cells.Item(24, 5).Value = 4922832168
$ws.Cells.Item(24, 6).Value = 98238057
$ws.Cells.Item(24, 7).Value = 0.01783
$ws.Cells.Item(25, 2).Value = "IMX"
$ws.Cells.Item(25, 3).Value = "Immutable"
$ws.Cells.Item(25, 4).Value = 3.2
$ws.Cells.Item(25, 5).Value = 4337330916
$ws.Cells.Item(25, 6).Value = 83119337
$ws.Cells.Item(25, 7).Value = 3.3533
$ws.Cells.Item(26, 2).Value = "TAO"
$ws.Cells.Item(26, 3).Value = "Bittensor"
$ws.Cells.Item(26, 4).Value = 649.87
$ws.Cells.Item(26, 5).Value = 4059177904
$ws.Cells.Item(26, 6).Value = 21971147
$ws.Cells.Item(26, 7).Value = 1.85878
$ws.Cells.Item(27, 2).Value = "ATOM"
$ws.Cells.Item(27, 3).Value = "Cosmos Hub"
$ws.Cells.Item(27, 4).Value = 10.39
$ws.Cells.Item(27, 5).Value = 3990499729
$ws.Cells.Item(27, 6).Value = 176653450
$ws.Cells.Item(27, 7).Value = 2.12098
$ws.Cells.Item(28, 2).Value = "LEO"
$ws.Cells.Item(28, 3).Value = "LEO Token"
$ws.Cells.Item(28, 4).Value = 4.15
$ws.Cells.Item(28, 5).Value = 3849487386
$ws.Cells.Item(28, 6).Value = 1119726
$ws.Cells.Item(28, 7).Value = 0.46465
$ws.Cells.Item(29, 2).Value = "ETC"
$ws.Cells.Item(29, 3).Value = "Ethereum Classic"
$ws.Cells.Item(29, 4).Value = 26.24
$ws.Cells.Item(29, 5).Value = 3756910136
$ws.Cells.Item(29, 6).Value = 139489865
$ws.Cells.Item(29, 7).Value = -0.4089
$ws.Cells.Item(30, 2).Value = "KAS"
$ws.Cells.Item(30, 3).Value = "Kaspa"
$ws.Cells.Item(30, 4).Value = 0.162431
$ws.Cells.Item(30, 5).Value = 3694535422
$ws.Cells.Item(30, 6).Value = 72986537
$ws.Cells.Item(30, 7).Value = -3.29102
$ws.Cells.Item(31, 2).Value = "STX"
$ws.Cells.Item(31, 3).Value = "Stacks"
$ws.Cells.Item(31, 4).Value = 2.51
$ws.Cells.Item(31, 5).Value = 3611508463
$ws.Cells.Item(31, 6).Value = 95586507
$ws.Cells.Item(31, 7).Value = 0.13227
$ws.Cells.Item(32, 2).Value = "OP"
$ws.Cells.Item(32, 3).Value = "Optimism"
$ws.Cells.Item(32, 4).Value = 3.72
$ws.Cells.Item(32, 5).Value = 3564175754
$ws.Cells.Item(32, 6).Value = 139695477
$ws.Cells.Item(32, 7).Value = 4.19683
$ws.Cells.Item(33, 2).Value = "NEAR"
$ws.Cells.Item(33, 3).Value = "NEAR Protocol"
$ws.Cells.Item(33, 4).Value = 3.43
$ws.Cells.Item(33, 5).Value = 3564001921
$ws.Cells.Item(33, 6).Value = 246251948
$ws.Cells.Item(33, 7).Value = 5.28741
$ws.Cells.Item(34, 2).Value = "APT"
$ws.Cells.Item(34, 3).Value = "Aptos"
$ws.Cells.Item(34, 4).Value = 9.699999999999999
$ws.Cells.Item(34, 5).Value = 3541783848
$ws.Cells.Item(34, 6).Value = 114642775
$ws.Cells.Item(34, 7).Value = 1.96908
$ws.Cells.Item(35, 2).Value = "INJ"
$ws.Cells.Item(35, 3).Value = "Injective"
$ws.Cells.Item(35, 4).Value = 38.95
$ws.Cells.Item(35, 5).Value = 3441400584
$ws.Cells.Item(35, 6).Value = 381088555
$ws.Cells.Item(35, 7).Value = 10.27923
$ws.Cells.Item(36, 2).Value = "XLM"
$ws.Cells.Item(36, 3).Value = "Stellar"
$ws.Cells.Item(36, 4).Value = 0.115572
$ws.Cells.Item(36, 5).Value = 3288692189
$ws.Cells.Item(36, 6).Value = 58416506
$ws.Cells.Item(36, 7).Value = 1.26964
$ws.Cells.Item(37, 2).Value = "VET"
$ws.Cells.Item(37, 3).Value = "VeChain"
$ws.Cells.Item(37, 4).Value = 0.04450379
$ws.Cells.Item(37, 5).Value = 3236308669
$ws.Cells.Item(37, 6).Value = 125530597
$ws.Cells.Item(37, 7).Value = 0.28737
$ws.Cells.Item(38, 2).Value = "FIL"
$ws.Cells.Item(38, 3).Value = "Filecoin"
$ws.Cells.Item(38, 4).Value = 6.17
$ws.Cells.Item(38, 5).Value = 3151154403
$ws.Cells.Item(38, 6).Value = 401962791
$ws.Cells.Item(38, 7).Value = -0.60083
$ws.Cells.Item(39, 2).Value = "OKB"
$ws.Cells.Item(39, 3).Value = "OKB"
$ws.Cells.Item(39, 4).Value = 52.42
$ws.Cells.Item(39, 5).Value = 3144932303
$ws.Cells.Item(39, 6).Value = 7777653
$ws.Cells.Item(39, 7).Value = 0.99581
$ws.Cells.Item(40, 2).Value = "TIA"
$ws.Cells.Item(40, 3).Value = "Celestia"
$ws.Cells.Item(40, 4).Value = 19.05
$ws.Cells.Item(40, 5).Value = 3143063039
$ws.Cells.Item(40, 6).Value = 125066530
$ws.Cells.Item(40, 7).Value = 1.74615
$ws.Cells.Item(41, 2).Value = "HBAR"
$ws.Cells.Item(41, 3).Value = "Hedera"
$ws.Cells.Item(41, 4).Value = 0.088579
$ws.Cells.Item(41, 5).Value = 2978859523
$ws.Cells.Item(41, 6).Value = 60830520
$ws.Cells.Item(41, 7).Value = 5.73682
$ws.Cells.Item(42, 2).Value = "FDUSD"
$ws.Cells.Item(42, 3).Value = "First Digital USD"
$ws.Cells.Item(42, 4).Value = 0.999685
$ws.Cells.Item(42, 5).Value = 2843180700
$ws.Cells.Item(42, 6).Value = 4052622152
$ws.Cells.Item(42, 7).Value = -0.03238
$ws.Cells.Item(43, 2).Value = "LDO"
$ws.Cells.Item(43, 3).Value = "Lido DAO"
$ws.Cells.Item(43, 4).Value = 3.17
$ws.Cells.Item(43, 5).Value = 2825234309
$ws.Cells.Item(43, 6).Value = 53234698
$ws.Cells.Item(43, 7).Value = 1.26788
$ws.Cells.Item(44, 2).Value = "ARB"
$ws.Cells.Item(44, 3).Value = "Arbitrum"
$ws.Cells.Item(44, 4).Value = 2.01
$ws.Cells.Item(44, 5).Value = 2561284754
$ws.Cells.Item(44, 6).Value = 231356886
$ws.Cells.Item(44, 7).Value = 3.0207
$ws.Cells.Item(45, 2).Value = "MNT"
$ws.Cells.Item(45, 3).Value = "Mantle"
$ws.Cells.Item(45, 4).Value = 0.777794
$ws.Cells.Item(45, 5).Value = 2505281240
$ws.Cells.Item(45, 6).Value = 80332201
$ws.Cells.Item(45, 7).Value = 2.92821
$ws.Cells.Item(46, 2).Value = "SEI"
$ws.Cells.Item(46, 3).Value = "Sei"
$ws.Cells.Item(46, 4).Value = 0.963441
$ws.Cells.Item(46, 5).Value = 2454237009
$ws.Cells.Item(46, 6).Value = 235534622
$ws.Cells.Item(46, 7).Value = 2.15667
$ws.Cells.Item(47, 2).Value = "CRO"
$ws.Cells.Item(47, 3).Value = "Cronos"
$ws.Cells.Item(47, 4).Value = 0.090457
$ws.Cells.Item(47, 5).Value = 2400524220
$ws.Cells.Item(47, 6).Value = 8120980
$ws.Cells.Item(47, 7).Value = 0.8381999999999999
$ws.Cells.Item(48, 2).Value = "XMR"
$ws.Cells.Item(48, 3).Value = "Monero"
$ws.Cells.Item(48, 4).Value = 120.64
$ws.Cells.Item(48, 5).Value = 2189011678
$ws.Cells.Item(48, 6).Value = 52566955
$ws.Cells.Item(48, 7).Value = 0.77472
$ws.Cells.Item(49, 2).Value = "RNDR"
$ws.Cells.Item(49, 3).Value = "Render"
$ws.Cells.Item(49, 4).Value = 5.63
$ws.Cells.Item(49, 5).Value = 2130294996
$ws.Cells.Item(49, 6).Value = 230671351
$ws.Cells.Item(49, 7).Value = 7.41635
$ws.Cells.Item(50, 2).Value = "SUI"
$ws.Cells.Item(50, 3).Value = "Sui"
$ws.Cells.Item(50, 4).Value = 1.79
$ws.Cells.Item(50, 5).Value = 2085640983
$ws.Cells.Item(50, 6).Value = 209089080
$ws.Cells.Item(50, 7).Value = 1.93564
$ws.Cells.Item(51, 2).Value = "GRT"
$ws.Cells.Item(51, 3).Value = "The Graph"
$ws.Cells.Item(51, 4).Value = 0.217561
$ws.Cells.Item(51, 5).Value = 2048881958
$ws.Cells.Item(51, 6).Value = 257755273
$ws.Cells.Item(51, 7).Value = 15.73712
